$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.212.93'
$ws.Range("E2").Value = '  +4.42%  '

$ws.Range("D3").Value = '1.707.33'
$ws.Range("E3").Value = '  +4.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.535'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.269'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0642'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0911'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.19%  '

$ws.Range("D13").Value = '1.947.93'
$ws.Range("E13").Value = '  +4.19%  '

$ws.Range("D14").Value = '1.708.68'
$ws.Range("E14").Value = '  +4.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.616'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.18%  '

$ws.Range("D18").Value = '31.193.35'
$ws.Range("E18").Value = '  +4.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '249.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.66%  '

$ws.Range("E21").Value = '  +2.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.72%  '

$ws.Range("E25").Value = '  -1.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.27%  '

$ws.Range("E28").Value = '  +3.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.77%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0504'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.25%  '

$ws.Range("E33").Value = '  +3.73%  '

$ws.Range("E34").Value = '  +6.06%  '

$ws.Range("D35").Value = '1.524.03'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.47%  '

$ws.Range("E37").Value = '  +1.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '83.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.01%  '

$ws.Range("E39").Value = '  +10.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0181'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.13%  '

$ws.Range("E41").Value = '  -2.49%  '

$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("E43").Value = '  +3.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.853'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("E46").Value = '  +3.08%  '

$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '52.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.73%  '

$ws.Range("D50").Value = '1.835.45'
$ws.Range("E50").Value = '  +3.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '94.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.97%  '
